$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet has a gap: row 4 is blank while row 5 holds the last
# data row ("kumash","kumash","admin"). Move that row up into row 4
# (note: bare ".Value" getter doesn't resolve in this host, so read
# with the explicit "Value()" call) and clear the now-vacated row 5.
$a5 = $ws.Range("A5").Value()
$b5 = $ws.Range("B5").Value()
$c5 = $ws.Range("C5").Value()

$ws.Range("A4").Value = $a5
$ws.Range("B4").Value = $b5
$ws.Range("C4").Value = $c5

$ws.Range("A5:C5").ClearContents()

# Match the saved selection: A4 active within A4:C4.
$ws.Range("A4:C4").Select()
